$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192707896232605
$ws.Range("B1").Value = 3.848259210586548
$ws.Range("C1").Value = 3.216675996780396
$ws.Range("D1").Value = 2.568390130996704
$ws.Range("E1").Value = 1.319860935211182
